$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (height 25.5): question + expected result (center+wrap) ---
$ws.Range("B17").Value2 = "Clicar no botão `"Fechar`" localizado ao lado do período para fechamento do semestre atual."
$ws.Range("B17").WrapText = $true
$ws.Range("C17").Value2 = "O status do semestre atual ficará igua a `"Fechado`"."
$ws.Range("C17").HorizontalAlignment = -4108
$ws.Range("C17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 25.5

# --- Row 18 ---
$ws.Range("B18").Value2 = "Clicar no botão `"Fechar trabalho`" ao lado de todos os trabalhos do semestre atual que estão com status de aberto."
$ws.Range("B18").WrapText = $true
$ws.Range("C18").Value2 = "Todos ostrabalhos referentes ao semestre atual serão fechados com sucesso."
$ws.Range("C18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 15.75

# --- Row 19 ---
$ws.Range("B19").Value2 = "Verificar no banco de dados o relacionamento entre os envolvidos do semestre atual"
$ws.Range("B19").WrapText = $true
$ws.Range("C19").Value2 = "Todos os relacionamentos deverão ser desfietos"
$ws.Range("C19").WrapText = $true
$ws.Range("C19").WrapText = $false
$ws.Rows.Item(19).RowHeight = 15.75

# --- Row 20 ---
$ws.Range("B20").Value2 = "Clicar na opção `"Adicionar novo semestre`"."
$ws.Range("B20").WrapText = $true
$ws.Range("C20").Value2 = "O sistema carrega a página informando os campos para preenchimento dos dados do novo semestre."
$ws.Range("C20").HorizontalAlignment = -4108
$ws.Range("C20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 15.75

# --- Row 21 ---
$ws.Range("B21").Value2 = "Preencher todos os campos obrigatórios e clicar em `"Cadastrar`" "
$ws.Range("B21").WrapText = $true
$ws.Range("C21").Value2 = "O sistema exibe uma mensagem informando que o semestre foi cadastrado com sucesso."
$ws.Range("C21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 15.75

# --- Row 22 ---
$ws.Range("B22").Value2 = "Verificar o status do novo semestre"
$ws.Range("B22").WrapText = $true
$ws.Range("C22").Value2 = "Após a criação do novo semestre, seu status deverá ser igual a `"atual`""
$ws.Range("C22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 15.75

# --- Row 23: "Fluxo de Exceção" section header (merged A23:C23) ---
$ws.Range("A23").Value2 = "Fluxo de Exceção"
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").Font.Name = "Arial"
$ws.Range("A23").Font.Size = 10
$ws.Range("A23").Interior.Color = 16308937
$ws.Range("A23").Interior.PatternColor = 16308937
$ws.Range("A23").HorizontalAlignment = -4108
$ws.Range("A23").WrapText = $true
$ws.Range("B23").Font.Name = "Arial"
$ws.Range("B23").Font.Size = 10
$ws.Range("C23").Font.Name = "Arial"
$ws.Range("C23").Font.Size = 10
$ws.Rows.Item(23).RowHeight = 15.75
$ws.Range("A23:C23").Merge()

# --- Row 24 ---
$ws.Range("B24").Value2 = "Preencher o campo `"Ano do semestre`" do novo semestre menor  que o ano do semestre que foi fechado"
$ws.Range("B24").WrapText = $true
$ws.Range("C24").Value2 = "Uma mensagem de erro deverá ser mostrada afirmando que não é possível criar um novo semestre com ano menor que o semestre anterior. O cadastro é cancelado."
$ws.Range("C24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 15.75

# --- Row 25 ---
$ws.Range("B25").Value2 = "Preencher os campos sem preencher os campos obrigatórios"
$ws.Range("B25").WrapText = $true
$ws.Range("C25").Value2 = "Uma mensagem de erro será mostrada informando que os campos obrigatórios deverão ser preenchido. Voltar pra tela de preenchimento dos campos"
$ws.Range("C25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 15.75

# --- Row 26 ---
$ws.Range("B26").Value2 = "Cadastrar um semestre já existente"
$ws.Range("B26").WrapText = $true
$ws.Range("C26").Value2 = "Uma mensagem de erro será mostrada informando que o semestre já existe. O cadastro é cancelado."
$ws.Range("C26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 15.75

# --- Rows 27-30: blank wrap-styled rows (B & C) ---
$ws.Range("B27").WrapText = $true
$ws.Range("C27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 15.75
$ws.Range("B28").WrapText = $true
$ws.Range("C28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 15.75
$ws.Range("B29").WrapText = $true
$ws.Range("C29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 15.75
$ws.Range("B30").WrapText = $true
$ws.Range("C30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 15.75

# --- Rows 31-32: blank wrap-styled rows (B only) ---
$ws.Range("B31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 15.75
$ws.Range("B32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 15.75

# --- Sheet view updates ---
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("B27").Select()

Write-Host "UC001 - Abrir e fechar semestre rows added"
